# Regenerate the K column (column G) values in the save_data sheet.
# Commit message: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (std-based strike-count), keyed by row number (row 1 is header).
$kValues = @{
    2 = 2
    3 = 0
    4 = 3
    5 = 0
    6 = 1
    7 = 0
    8 = 2
    9 = 1
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 3
    26 = 2
    27 = 0
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 2
    37 = 1
    38 = 0
    39 = 2
    40 = 2
    41 = 1
    42 = 0
    43 = 0
    44 = 1
    45 = 1
    46 = 1
    47 = 2
    48 = 2
    49 = 3
    50 = 2
    51 = 1
    52 = 2
    53 = 2
    54 = 0
    55 = 0
    56 = 1
    57 = 1
    58 = 3
    59 = 1
    60 = 1
    61 = 1
    62 = 0
    63 = 0
    64 = 2
    65 = 0
    66 = 2
    67 = 0
    68 = 1
    69 = 1
    70 = 1
    71 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
